$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column (match styling of neighboring header cell)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values, row 2 through 26
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 1
}

foreach ($row in 2..26) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
